# Add a "dSoH (ppm)" column (F/G) to the scenario comparison table, and
# correct the G2V "Trading rev." / "Adj. rev." figures from -2.8 to -2.5.
#
# Table layout (row 1 = generic numeric-string column ids, row 2 = real
# headers, rows 3-4 = G2V / V2G data):
#
#            A       B          C             D         E               F            G
#   1        0       1          2             3         4               5
#   2        0.0     Scenario   Trading rev.  HW cost   House el cost   Adj. rev.    dSoH (ppm)
#   3        1.0     G2V        -2.5          -0        0               -2.5         50
#   4        2.0     V2G        -1.9          -0.6      0               -2.5         50.1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Several of the new/changed values look like numbers ("5", "-2.5", "50",
# "50.1") but must be stored as literal text (matching every other cell in
# this table). Mark those cells as Text *before* writing to them so Excel
# doesn't silently convert them to numeric cells.
$textCells = @("F1", "C3", "F3", "G3", "G4")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# New column header on the index row.
$ws.Range("F1").Value = "5"

# New "dSoH (ppm)" column header (plain text, no coercion risk).
$ws.Range("G2").Value = "dSoH (ppm)"

# G2V row: trading revenue / adjusted revenue corrected, dSoH added.
$ws.Range("C3").Value = "-2.5"
$ws.Range("F3").Value = "-2.5"
$ws.Range("G3").Value = "50"

# V2G row: new dSoH value.
$ws.Range("G4").Value = "50.1"

# Match the new column's width to the author's resize (bestFit-style).
# (Excel's character-based ColumnWidth snaps to whole pixels, so 7.74 is
# the input that lands closest to the author's stored width of ~8.457.)
$ws.Columns.Item(6).ColumnWidth = 7.74
